$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.15"
$ws.Range("E2").Value = "'4.64%"
$ws.Range("D3").Value = "'26.89"
$ws.Range("E3").Value = "'0.74%"
$ws.Range("E4").Value = "'4.55%"
$ws.Range("D5").Value = "'0.06394"
$ws.Range("E5").Value = "'5.14%"
$ws.Range("D6").Value = "'6.998"
$ws.Range("E6").Value = "'3.82%"
$ws.Range("D7").Value = "'3.357"
$ws.Range("E7").Value = "'5.93%"
$ws.Range("D8").Value = "'0.8875"
$ws.Range("E8").Value = "'4.30%"
$ws.Range("D9").Value = "'1.174"
$ws.Range("E9").Value = "'29.55%"
$ws.Range("E10").Value = "'5.91%"
$ws.Range("D11").Value = "'0.05265"
$ws.Range("E11").Value = "'5.35%"
$ws.Range("D12").Value = "'0.07418"
$ws.Range("E12").Value = "'4.62%"
$ws.Range("D13").Value = "'0.03150"
$ws.Range("E13").Value = "'-1.19%"
$ws.Range("D14").Value = "'0.09066"
$ws.Range("E14").Value = "'0.60%"
$ws.Range("D15").Value = "'0.001561"
$ws.Range("E15").Value = "'1.55%"
$ws.Range("D16").Value = "'0.0006347"
$ws.Range("E16").Value = "'5.01%"
$ws.Range("D17").Value = "'0.006041"
$ws.Range("E17").Value = "'0.42%"
$ws.Range("D18").Value = "'3.492"
$ws.Range("E18").Value = "'0.91%"
$ws.Range("D19").Value = "'2.280"
$ws.Range("E19").Value = "'0.72%"
$ws.Range("D20").Value = "'0.3156"
$ws.Range("E20").Value = "'2.19%"
$ws.Range("D21").Value = "'0.1332"
$ws.Range("E21").Value = "'2.51%"
$ws.Range("D22").Value = "'3.925"
$ws.Range("E22").Value = "'-3.69%"
$ws.Range("D23").Value = "'0.04355"
$ws.Range("E23").Value = "'2.97%"
$ws.Range("E24").Value = "'0.02%"
$ws.Range("E25").Value = "'-11.03%"
$ws.Range("D27").Value = "'0.0001701"
$ws.Range("E27").Value = "'1.15%"
$ws.Range("D40").Value = "'0.04070"
$ws.Range("E40").Value = "'3.94%"
$ws.Range("D41").Value = "'0.006646"
$ws.Range("E41").Value = "'59.08%"
$ws.Range("D42").Value = "'0.1174"
$ws.Range("E42").Value = "'5.45%"
$ws.Range("D43").Value = "'0.002364"
$ws.Range("E43").Value = "'11.97%"
$ws.Range("D44").Value = "'0.01286"
$ws.Range("E44").Value = "'1.84%"
$ws.Range("D45").Value = "'0.00005271"
$ws.Range("E45").Value = "'3.33%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E47").Value = "'1,625.68%"
$ws.Range("D48").Value = "'0.02124"
$ws.Range("E48").Value = "'-13.24%"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'-0.11%"
